$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('展览')
$ws.Rows.Item(2).Delete()
$ws.Range("B2:B9").NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = '2024-05-19'
$ws.Cells.Item(2, 3).Value = '南宁·原x穹x崩only'
$ws.Cells.Item(2, 4).Value = '明秀东路157号 利泰国际大酒店'
$ws.Cells.Item(2, 5).Value = '2024.05.19 10:00-05.19 17:00'
$ws.Cells.Item(2, 6).Value = 266
$ws.Cells.Item(2, 7).Value = 35
$ws.Cells.Item(2, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=83070'
$ws.Cells.Item(2, 9).Value = '//i0.hdslb.com/bfs/openplatform/202403/I8tScigE1710918412731.jpeg'
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = '2024-05-25'
$ws.Cells.Item(3, 3).Value = '南宁·第五人格Only1.0'
$ws.Cells.Item(3, 4).Value = '新阳路227号南宁第三人民医院旁新秀佳园对面 卡尔顿东方银龙酒店'
$ws.Cells.Item(3, 5).Value = '2024.05.25 10:00-05.25 17:30'
$ws.Cells.Item(3, 6).Value = 156
$ws.Cells.Item(3, 7).Value = 68
$ws.Cells.Item(3, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84954'
$ws.Cells.Item(3, 9).Value = '//i0.hdslb.com/bfs/openplatform/202404/w5iZT4wE1714189905443.jpeg'
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = '2024-06-09'
$ws.Cells.Item(4, 3).Value = '南宁·AP动漫游戏嘉年华'
$ws.Cells.Item(4, 4).Value = '三塘南路与长虹东路交叉路口往北约50米 广西农业会展中心'
$ws.Cells.Item(4, 5).Value = '2024.06.09 09:00-06.10 17:00'
$ws.Cells.Item(4, 6).Value = 1850
$ws.Cells.Item(4, 7).Value = 60
$ws.Cells.Item(4, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84793'
$ws.Cells.Item(4, 9).Value = '//i0.hdslb.com/bfs/openplatform/202404/lNO8AZwt1712818829933.jpeg'
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = '2024-06-09'
$ws.Cells.Item(5, 3).Value = '南宁·布谷鸟动漫展4th'
$ws.Cells.Item(5, 4).Value = '亭洪路45号 百益上河城'
$ws.Cells.Item(5, 5).Value = '2024.06.09 10:00-06.10 17:00'
$ws.Cells.Item(5, 6).Value = 1557
$ws.Cells.Item(5, 7).Value = 50
$ws.Cells.Item(5, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82241'
$ws.Cells.Item(5, 9).Value = '//i2.hdslb.com/bfs/openplatform/202403/uzZqZov91709281147333.jpeg'
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = '2024-06-09'
$ws.Cells.Item(6, 3).Value = '南宁·恋与深空only'
$ws.Cells.Item(6, 4).Value = '新阳路227号南宁第三人民医院旁新秀佳园对面 卡尔顿东方银龙酒店'
$ws.Cells.Item(6, 5).Value = '2024.06.09 10:00-06.09 17:00'
$ws.Cells.Item(6, 6).Value = 282
$ws.Cells.Item(6, 7).Value = 50
$ws.Cells.Item(6, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84444'
$ws.Cells.Item(6, 9).Value = '//i2.hdslb.com/bfs/openplatform/202404/6ZVHU1F91713340880421.jpeg'
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = '2024-07-06'
$ws.Cells.Item(7, 3).Value = '南宁·小蜜蜂动漫嘉年华2.0'
$ws.Cells.Item(7, 4).Value = '亭洪路45号 百益上河城'
$ws.Cells.Item(7, 5).Value = '2024.07.06 10:00-07.06 17:00'
$ws.Cells.Item(7, 6).Value = 70
$ws.Cells.Item(7, 7).Value = 50
$ws.Cells.Item(7, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84925'
$ws.Cells.Item(7, 9).Value = '//i2.hdslb.com/bfs/openplatform/202404/YjFyyYq51713508727131.jpeg'
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = '2024-07-20'
$ws.Cells.Item(8, 3).Value = '南宁·AB动漫游戏嘉年华'
$ws.Cells.Item(8, 4).Value = '三塘南路与长虹东路交叉路口往北约50米 广西农业会展中心'
$ws.Cells.Item(8, 5).Value = '2024.07.20 09:30-07.21 17:00'
$ws.Cells.Item(8, 6).Value = 538
$ws.Cells.Item(8, 7).Value = 60
$ws.Cells.Item(8, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84862'
$ws.Cells.Item(8, 9).Value = '//i1.hdslb.com/bfs/openplatform/202404/eglavDeZ1714036487217.jpeg'
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = '2024-07-20'
$ws.Cells.Item(9, 3).Value = '横州·第二届海棠动漫游戏嘉年华'
$ws.Cells.Item(9, 4).Value = '茉莉花大道 横州国际大酒店'
$ws.Cells.Item(9, 5).Value = '2024.07.20 09:30-07.20 17:00'
$ws.Cells.Item(9, 6).Value = 130
$ws.Cells.Item(9, 7).Value = 30
$ws.Cells.Item(9, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84799'
$ws.Cells.Item(9, 9).Value = '//i2.hdslb.com/bfs/openplatform/202404/r50S2ttT1713869164413.jpeg'

$ws = $wb.Worksheets.Item('全部类型')
$ws.Rows.Item(2).Delete()
$ws.Range("B2:B10").NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = '2024-05-19'
$ws.Cells.Item(2, 3).Value = '南宁·原x穹x崩only'
$ws.Cells.Item(2, 4).Value = '明秀东路157号 利泰国际大酒店'
$ws.Cells.Item(2, 5).Value = '2024.05.19 10:00-05.19 17:00'
$ws.Cells.Item(2, 6).Value = 266
$ws.Cells.Item(2, 7).Value = 35
$ws.Cells.Item(2, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=83070'
$ws.Cells.Item(2, 9).Value = '//i0.hdslb.com/bfs/openplatform/202403/I8tScigE1710918412731.jpeg'
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = '2024-05-25'
$ws.Cells.Item(3, 3).Value = '南宁·第五人格Only1.0'
$ws.Cells.Item(3, 4).Value = '新阳路227号南宁第三人民医院旁新秀佳园对面 卡尔顿东方银龙酒店'
$ws.Cells.Item(3, 5).Value = '2024.05.25 10:00-05.25 17:30'
$ws.Cells.Item(3, 6).Value = 156
$ws.Cells.Item(3, 7).Value = 68
$ws.Cells.Item(3, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84954'
$ws.Cells.Item(3, 9).Value = '//i0.hdslb.com/bfs/openplatform/202404/w5iZT4wE1714189905443.jpeg'
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = '2024-06-09'
$ws.Cells.Item(4, 3).Value = '南宁·AP动漫游戏嘉年华'
$ws.Cells.Item(4, 4).Value = '三塘南路与长虹东路交叉路口往北约50米 广西农业会展中心'
$ws.Cells.Item(4, 5).Value = '2024.06.09 09:00-06.10 17:00'
$ws.Cells.Item(4, 6).Value = 1850
$ws.Cells.Item(4, 7).Value = 60
$ws.Cells.Item(4, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84793'
$ws.Cells.Item(4, 9).Value = '//i0.hdslb.com/bfs/openplatform/202404/lNO8AZwt1712818829933.jpeg'
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = '2024-06-09'
$ws.Cells.Item(5, 3).Value = '南宁·布谷鸟动漫展4th'
$ws.Cells.Item(5, 4).Value = '亭洪路45号 百益上河城'
$ws.Cells.Item(5, 5).Value = '2024.06.09 10:00-06.10 17:00'
$ws.Cells.Item(5, 6).Value = 1557
$ws.Cells.Item(5, 7).Value = 50
$ws.Cells.Item(5, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82241'
$ws.Cells.Item(5, 9).Value = '//i2.hdslb.com/bfs/openplatform/202403/uzZqZov91709281147333.jpeg'
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = '2024-06-09'
$ws.Cells.Item(6, 3).Value = '南宁·恋与深空only'
$ws.Cells.Item(6, 4).Value = '新阳路227号南宁第三人民医院旁新秀佳园对面 卡尔顿东方银龙酒店'
$ws.Cells.Item(6, 5).Value = '2024.06.09 10:00-06.09 17:00'
$ws.Cells.Item(6, 6).Value = 282
$ws.Cells.Item(6, 7).Value = 50
$ws.Cells.Item(6, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84444'
$ws.Cells.Item(6, 9).Value = '//i2.hdslb.com/bfs/openplatform/202404/6ZVHU1F91713340880421.jpeg'
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = '2024-06-22'
$ws.Cells.Item(7, 3).Value = '南宁·浪漫古典·百年经典世界名曲音乐会'
$ws.Cells.Item(7, 4).Value = '广西壮族自治区南宁市良庆区龙堤路25号  广西文化艺术中心-音乐厅'
$ws.Cells.Item(7, 5).Value = '2024.06.22 20:00-06.22 21:30'
$ws.Cells.Item(7, 6).Value = 19
$ws.Cells.Item(7, 7).Value = 50
$ws.Cells.Item(7, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=83959'
$ws.Cells.Item(7, 9).Value = '//i1.hdslb.com/bfs/openplatform/202404/H0f8U7no1712041461015.jpeg'
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = '2024-07-06'
$ws.Cells.Item(8, 3).Value = '南宁·小蜜蜂动漫嘉年华2.0'
$ws.Cells.Item(8, 4).Value = '亭洪路45号 百益上河城'
$ws.Cells.Item(8, 5).Value = '2024.07.06 10:00-07.06 17:00'
$ws.Cells.Item(8, 6).Value = 70
$ws.Cells.Item(8, 7).Value = 50
$ws.Cells.Item(8, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84925'
$ws.Cells.Item(8, 9).Value = '//i2.hdslb.com/bfs/openplatform/202404/YjFyyYq51713508727131.jpeg'
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = '2024-07-20'
$ws.Cells.Item(9, 3).Value = '南宁·AB动漫游戏嘉年华'
$ws.Cells.Item(9, 4).Value = '三塘南路与长虹东路交叉路口往北约50米 广西农业会展中心'
$ws.Cells.Item(9, 5).Value = '2024.07.20 09:30-07.21 17:00'
$ws.Cells.Item(9, 6).Value = 538
$ws.Cells.Item(9, 7).Value = 60
$ws.Cells.Item(9, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84862'
$ws.Cells.Item(9, 9).Value = '//i1.hdslb.com/bfs/openplatform/202404/eglavDeZ1714036487217.jpeg'
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = '2024-07-20'
$ws.Cells.Item(10, 3).Value = '横州·第二届海棠动漫游戏嘉年华'
$ws.Cells.Item(10, 4).Value = '茉莉花大道 横州国际大酒店'
$ws.Cells.Item(10, 5).Value = '2024.07.20 09:30-07.20 17:00'
$ws.Cells.Item(10, 6).Value = 130
$ws.Cells.Item(10, 7).Value = 30
$ws.Cells.Item(10, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84799'
$ws.Cells.Item(10, 9).Value = '//i2.hdslb.com/bfs/openplatform/202404/r50S2ttT1713869164413.jpeg'

Write-Output "done"
